$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D used to hold "panel" - replace it with the street address and
# shift the remaining contact/location fields over to make room for the
# new city/state/zip/coordinate columns.
$ws.Range("D1").Value = "2213 E 2100 S"
$ws.Range("E1").Value = "Salt Lake City"
$ws.Range("F1").Value = "Utah"
$ws.Range("G1").Value = 84102
$ws.Range("H1").Value = 40.72589
$ws.Range("I1").Value = -111.82782
$ws.Range("J1").Value = $false

# Give C2 the built-in "Hyperlink" cell style (empty placeholder cell for a
# future link) without leaving a real hyperlink registered on the sheet.
$link = $ws.Hyperlinks.Add($ws.Range("C2"), "", "", "", "")
$link.Delete()

# Mirror the saved selection state (row 2 highlighted A2:F2, anchored at F2).
$ws.Range("A2:F2").Select()
